$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The population column previously read "basisschoolleerlingen van groep 8"
# for the c11_* outcome rows (5-18); rename it to "leerlingen van groep 8".
$ws.Range("D5:D18").Value = "leerlingen van groep 8"

# Reflect the cell that was left selected in the sheet.
$ws.Range("D20").Select()
